$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")
$labels  = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $labels[$i] + "_FV2410"
    $ws.Range($newCols[$i] + "1").Value = $labels[$i] + "_FV2504"
}

# --- 2. Turn the used range into an Excel Table (ListObject) named "Table1" ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U62"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
